$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-09-04T15:07:43"
$ws.Range("S4").Value = 44.34
$ws.Range("T4").Value = 40.44
$ws.Range("U4").Value = 46.15
$ws.Range("V4").Value = 44.96
$ws.Range("W4").Value = 36.71
$ws.Range("X4").Value = 36.14
$ws.Range("Y4").Value = 31.77
$ws.Range("Z4").Value = 30.33
$ws.Range("S6").Value = -1.77
$ws.Range("T6").Value = -1.62
$ws.Range("U6").Value = -1.75
$ws.Range("V6").Value = -1.84
$ws.Range("W6").Value = -1.28
$ws.Range("Y6").Value = -0.79
$ws.Range("Z6").Value = -0.76
$ws.Range("S9").Value = 44.09
$ws.Range("T9").Value = 40.67
$ws.Range("U9").Value = 46.55
$ws.Range("V9").Value = 45.31
$ws.Range("W9").Value = 37.85
$ws.Range("Y9").Value = 32.79
$ws.Range("Z9").Value = 31.24
$ws.Range("S11").Value = -2.03
$ws.Range("T11").Value = -1.38
$ws.Range("U11").Value = -1.35
$ws.Range("V11").Value = -1.5
$ws.Range("W11").Value = -0.15
$ws.Range("X11").Value = -0.66
$ws.Range("Y11").Value = 0.23
$ws.Range("Z11").Value = 0.16
$ws.Range("S14").Value = 44.09
$ws.Range("T14").Value = 40.71
$ws.Range("U14").Value = 46.55
$ws.Range("V14").Value = 45.31
$ws.Range("W14").Value = 37.85
$ws.Range("Y14").Value = 32.79
$ws.Range("Z14").Value = 31.24
$ws.Range("S16").Value = -2.03
$ws.Range("T16").Value = -1.34
$ws.Range("U16").Value = -1.35
$ws.Range("V16").Value = -1.5
$ws.Range("X16").Value = -0.66
$ws.Range("Y16").Value = 0.23
$ws.Range("Z16").Value = 0.16
$ws.Range("S19").Value = 44.13
$ws.Range("T19").Value = 40.28
$ws.Range("U19").Value = 45.97
$ws.Range("V19").Value = 44.74
$ws.Range("W19").Value = 36.64
$ws.Range("Y19").Value = 31.77
$ws.Range("Z19").Value = 30.36
$ws.Range("S21").Value = -1.99
$ws.Range("T21").Value = -1.77
$ws.Range("U21").Value = -1.93
$ws.Range("V21").Value = -2.06
$ws.Range("W21").Value = -1.36
$ws.Range("X21").Value = -1.37
$ws.Range("Y21").Value = -0.79
$ws.Range("Z21").Value = -0.73
$ws.Range("S24").Value = 44.13
$ws.Range("T24").Value = 40.28
$ws.Range("U24").Value = 45.97
$ws.Range("V24").Value = 44.74
$ws.Range("W24").Value = 36.64
$ws.Range("Y24").Value = 31.77
$ws.Range("Z24").Value = 30.36
$ws.Range("S26").Value = -1.99
$ws.Range("T26").Value = -1.77
$ws.Range("U26").Value = -1.93
$ws.Range("V26").Value = -2.06
$ws.Range("W26").Value = -1.36
$ws.Range("X26").Value = -1.37
$ws.Range("Y26").Value = -0.79
$ws.Range("Z26").Value = -0.73
$ws.Range("S29").Value = 43.88
$ws.Range("T29").Value = 40.09
$ws.Range("U29").Value = 45.8
$ws.Range("V29").Value = 44.53
$ws.Range("W29").Value = 36.57
$ws.Range("X29").Value = 35.9
$ws.Range("Y29").Value = 31.77
$ws.Range("Z29").Value = 30.39
$ws.Range("S31").Value = -2.24
$ws.Range("T31").Value = -1.96
$ws.Range("U31").Value = -2.11
$ws.Range("V31").Value = -2.27
$ws.Range("W31").Value = -1.43
$ws.Range("Y31").Value = -0.79
$ws.Range("Z31").Value = -0.7
$ws.Range("S34").Value = 44.09
$ws.Range("T34").Value = 40.95
$ws.Range("U34").Value = 46.92
$ws.Range("V34").Value = 45.7
$ws.Range("W34").Value = 38.5
$ws.Range("Y34").Value = 33.43
$ws.Range("Z34").Value = 31.78
$ws.Range("S36").Value = -2.03
$ws.Range("T36").Value = -1.11
$ws.Range("U36").Value = -0.99
$ws.Range("V36").Value = -1.1
$ws.Range("X36").Value = -0.22
$ws.Range("Y36").Value = 0.87
$ws.Range("Z36").Value = 0.7
$ws.Range("S39").Value = 44.34
$ws.Range("T39").Value = 40.44
$ws.Range("U39").Value = 46.15
$ws.Range("V39").Value = 44.96
$ws.Range("W39").Value = 36.71
$ws.Range("X39").Value = 36.14
$ws.Range("Y39").Value = 31.77
$ws.Range("Z39").Value = 30.33
$ws.Range("S41").Value = -1.77
$ws.Range("T41").Value = -1.62
$ws.Range("U41").Value = -1.75
$ws.Range("V41").Value = -1.84
$ws.Range("W41").Value = -1.28
$ws.Range("Y41").Value = -0.79
$ws.Range("Z41").Value = -0.76
$ws.Range("S44").Value = 46.46
$ws.Range("T44").Value = 42.31
$ws.Range("U44").Value = 48.22
$ws.Range("V44").Value = 47.08
$ws.Range("W44").Value = 38.26
$ws.Range("X44").Value = 37.76
$ws.Range("Y44").Value = 32.97
$ws.Range("Z44").Value = 31.39
$ws.Range("T46").Value = 0.25
$ws.Range("U46").Value = 0.32
$ws.Range("V46").Value = 0.28
$ws.Range("W46").Value = 0.26
$ws.Range("Y46").Value = 0.41
$ws.Range("Z46").Value = 0.3
$ws.Range("S49").Value = 45.66
$ws.Range("T49").Value = 41.39
$ws.Range("U49").Value = 49.59
$ws.Range("V49").Value = 48.45
$ws.Range("W49").Value = 39.71
$ws.Range("X49").Value = 38.8
$ws.Range("Y49").Value = 33.68
$ws.Range("Z49").Value = 31.98
$ws.Range("S51").Value = -0.46
$ws.Range("T51").Value = -0.66
$ws.Range("U51").Value = 1.69
$ws.Range("V51").Value = 1.65
$ws.Range("W51").Value = 1.71
$ws.Range("X51").Value = 1.44
$ws.Range("Y51").Value = 1.11
$ws.Range("Z51").Value = 0.9
$ws.Range("S54").Value = 45.89
$ws.Range("T54").Value = 42.1
$ws.Range("U54").Value = 48.05
$ws.Range("V54").Value = 47.85
$ws.Range("W54").Value = 38.97
$ws.Range("X54").Value = 38.17
$ws.Range("Y54").Value = 33.16
$ws.Range("Z54").Value = 31.53
$ws.Range("T56").Value = 0.04
$ws.Range("U56").Value = 0.14
$ws.Range("V56").Value = 1.05
$ws.Range("W56").Value = 0.97
$ws.Range("Y56").Value = 0.6
$ws.Range("S59").Value = 48.24
$ws.Range("T59").Value = 43.86
$ws.Range("U59").Value = 50
$ws.Range("V59").Value = 48.85
$ws.Range("W59").Value = 39.62
$ws.Range("X59").Value = 39.13
$ws.Range("Y59").Value = 34.1
$ws.Range("Z59").Value = 32.45
$ws.Range("S61").Value = 2.12
$ws.Range("T61").Value = 1.8
$ws.Range("U61").Value = 2.1
$ws.Range("V61").Value = 2.05
$ws.Range("W61").Value = 1.62
$ws.Range("Y61").Value = 1.53
$ws.Range("Z61").Value = 1.36
$ws.Range("S64").Value = 48.75
$ws.Range("T64").Value = 44.32
$ws.Range("U64").Value = 50.53
$ws.Range("V64").Value = 49.37
$ws.Range("W64").Value = 40.04
$ws.Range("Y64").Value = 34.35
$ws.Range("Z64").Value = 32.69
$ws.Range("S66").Value = 2.63
$ws.Range("T66").Value = 2.26
$ws.Range("U66").Value = 2.63
$ws.Range("V66").Value = 2.57
$ws.Range("W66").Value = 2.04
$ws.Range("X66").Value = 2.09
$ws.Range("Y66").Value = 1.79
$ws.Range("Z66").Value = 1.6
$ws.Range("S69").Value = 49.27
$ws.Range("T69").Value = 44.84
$ws.Range("U69").Value = 51.07
$ws.Range("V69").Value = 49.89
$ws.Range("W69").Value = 40.51
$ws.Range("Y69").Value = 34.83
$ws.Range("Z69").Value = 33.17
$ws.Range("S71").Value = 3.15
$ws.Range("T71").Value = 2.78
$ws.Range("U71").Value = 3.17
$ws.Range("V71").Value = 3.09
$ws.Range("W71").Value = 2.51
$ws.Range("X71").Value = 2.6
$ws.Range("Y71").Value = 2.26
$ws.Range("Z71").Value = 2.09
$ws.Range("S74").Value = 47.79
$ws.Range("T74").Value = 43.49
$ws.Range("U74").Value = 49.49
$ws.Range("V74").Value = 48.4
$ws.Range("W74").Value = 39.38
$ws.Range("X74").Value = 38.8
$ws.Range("Y74").Value = 33.78
$ws.Range("Z74").Value = 32.18
$ws.Range("S76").Value = 1.67
$ws.Range("T76").Value = 1.44
$ws.Range("U76").Value = 1.58
$ws.Range("V76").Value = 1.6
$ws.Range("W76").Value = 1.38
$ws.Range("X76").Value = 1.44
$ws.Range("Y76").Value = 1.22
$ws.Range("S79").Value = 48.07
$ws.Range("T79").Value = 43.67
$ws.Range("U79").Value = 49.71
$ws.Range("V79").Value = 48.69
$ws.Range("W79").Value = 39.56
$ws.Range("X79").Value = 39.02
$ws.Range("Y79").Value = 33.99
$ws.Range("Z79").Value = 32.37
$ws.Range("S81").Value = 1.96
$ws.Range("T81").Value = 1.62
$ws.Range("U81").Value = 1.8
$ws.Range("V81").Value = 1.88
$ws.Range("W81").Value = 1.57
$ws.Range("X81").Value = 1.66
$ws.Range("Y81").Value = 1.43
$ws.Range("Z81").Value = 1.28
$ws.Range("S84").Value = 44.47
$ws.Range("T84").Value = 40.87
$ws.Range("U84").Value = 46.87
$ws.Range("V84").Value = 47.95
$ws.Range("W84").Value = 39.09
$ws.Range("X84").Value = 38.09
$ws.Range("Y84").Value = 33.13
$ws.Range("Z84").Value = 31.59
$ws.Range("S86").Value = -1.65
$ws.Range("T86").Value = -1.19
$ws.Range("U86").Value = -1.03
$ws.Range("V86").Value = 1.15
$ws.Range("W86").Value = 1.09
$ws.Range("Y86").Value = 0.5600000000000001
$ws.Range("Z86").Value = 0.51
$ws.Range("S89").Value = 43.88
$ws.Range("T89").Value = 40.09
$ws.Range("U89").Value = 45.8
$ws.Range("V89").Value = 44.53
$ws.Range("W89").Value = 36.57
$ws.Range("X89").Value = 35.9
$ws.Range("Y89").Value = 31.77
$ws.Range("Z89").Value = 30.39
$ws.Range("S91").Value = -2.24
$ws.Range("T91").Value = -1.96
$ws.Range("U91").Value = -2.11
$ws.Range("V91").Value = -2.27
$ws.Range("W91").Value = -1.43
$ws.Range("Y91").Value = -0.79
$ws.Range("Z91").Value = -0.7

Write-Output "Applied all cell updates"